# Update column F (dSF) values for specific rows, per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -2
$ws.Range("F15").Value = -5
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -6
$ws.Range("F23").Value = -2
$ws.Range("F25").Value = -4
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = -6
$ws.Range("F32").Value = -15
$ws.Range("F33").Value = -2
